$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark from its original location
#    (right after "...AdventureWorks case study").
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Update row "7" of the requirements table: the data requirement
#    text and the technology text change.
# ------------------------------------------------------------------
$t = $d.Tables.Item(1)

$t.Cell(8, 2).Range.Text = "The customer services department want to help their agents to identify fraudulent call for support"
$t.Cell(8, 3).Range.Text = "Stream Analytics / Event Hubs"

# Re-create the "_GoBack" bookmark right after the new technology text.
$techRange = $t.Cell(8, 3).Range
$techRange.MoveEnd(1, -1) | Out-Null
$bmRange = $techRange.Duplicate
$bmRange.Collapse(0) | Out-Null
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# ------------------------------------------------------------------
# 3. Insert a new row (becomes row "8") directly below row "7",
#    re-using the text that used to live in row "7" before the edit.
# ------------------------------------------------------------------
$beforeRow = $t.Rows.Item(9)
$t.Rows.Add($beforeRow) | Out-Null

$t.Cell(9, 1).Range.Text = "8"
$t.Cell(9, 2).Range.Text = "First party and 3rd party applications have access to the information of the bicycle computer"
$t.Cell(9, 3).Range.Text = "Stream Analytics / IoT Hubs"

# ------------------------------------------------------------------
# 4. Renumber the old row "8" (Bike telemetry ...) to "9".
# ------------------------------------------------------------------
$t.Cell(10, 1).Range.Text = "9"
